# Replace the "car/dog" working set of cue sequences (rows 2-33, columns
# B=image, C=word, D=category) with the new "dog/house" working set.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "dog/dog004.jpg"
$ws.Range("C2").Value = "klappen"
$ws.Range("D2").Value = "dog"

$ws.Range("B3").Value = "dog/dog017.jpg"
$ws.Range("C3").Value = "biegen"
$ws.Range("D3").Value = "dog"

$ws.Range("B4").Value = "dog/dog000.jpg"
$ws.Range("C4").Value = "lehnen"
$ws.Range("D4").Value = "dog"

$ws.Range("B5").Value = "dog/dog007.jpg"
$ws.Range("C5").Value = "loben"
$ws.Range("D5").Value = "dog"

$ws.Range("B6").Value = "house/house011.jpg"
$ws.Range("C6").Value = "dauern"
$ws.Range("D6").Value = "house"

$ws.Range("B7").Value = "dog/dog010.jpg"
$ws.Range("C7").Value = "strahlen"
$ws.Range("D7").Value = "dog"

$ws.Range("B8").Value = "house/house031.jpg"
$ws.Range("C8").Value = "schmecken"
$ws.Range("D8").Value = "house"

$ws.Range("B9").Value = "house/house009.jpg"
$ws.Range("C9").Value = "runden"
$ws.Range("D9").Value = "house"

$ws.Range("B10").Value = "dog/dog016.jpg"
$ws.Range("C10").Value = "posten"
$ws.Range("D10").Value = "dog"

$ws.Range("B11").Value = "house/house008.jpg"
$ws.Range("C11").Value = "enden"
$ws.Range("D11").Value = "house"

$ws.Range("B12").Value = "dog/dog013.jpg"
$ws.Range("C12").Value = "formen"
$ws.Range("D12").Value = "dog"

$ws.Range("B13").Value = "house/house025.jpg"
$ws.Range("C13").Value = "fliehen"
$ws.Range("D13").Value = "house"

$ws.Range("B14").Value = "dog/dog031.jpg"
$ws.Range("C14").Value = "antun"
$ws.Range("D14").Value = "dog"

$ws.Range("B15").Value = "house/house030.jpg"
$ws.Range("C15").Value = "mieten"
$ws.Range("D15").Value = "house"

$ws.Range("B16").Value = "house/house023.jpg"
$ws.Range("C16").Value = "nehmen"
$ws.Range("D16").Value = "house"

$ws.Range("B17").Value = "house/house024.jpg"
$ws.Range("C17").Value = "füllen"
$ws.Range("D17").Value = "house"

$ws.Range("B18").Value = "dog/dog030.jpg"
$ws.Range("C18").Value = "fliegen"
$ws.Range("D18").Value = "dog"

$ws.Range("B19").Value = "house/house019.jpg"
$ws.Range("C19").Value = "wenden"
$ws.Range("D19").Value = "house"

$ws.Range("B20").Value = "house/house018.jpg"
$ws.Range("C20").Value = "drohen"
$ws.Range("D20").Value = "house"

$ws.Range("B21").Value = "house/house022.jpg"
$ws.Range("C21").Value = "füttern"
$ws.Range("D21").Value = "house"

$ws.Range("B22").Value = "house/house013.jpg"
$ws.Range("C22").Value = "hauen"
$ws.Range("D22").Value = "house"

$ws.Range("B23").Value = "house/house003.jpg"
$ws.Range("C23").Value = "sondern"
$ws.Range("D23").Value = "house"

$ws.Range("B24").Value = "house/house020.jpg"
$ws.Range("C24").Value = "rasen"
$ws.Range("D24").Value = "house"

$ws.Range("B25").Value = "house/house015.jpg"
$ws.Range("C25").Value = "drehen"
$ws.Range("D25").Value = "house"

$ws.Range("B26").Value = "dog/dog015.jpg"
$ws.Range("C26").Value = "schätzen"
$ws.Range("D26").Value = "dog"

$ws.Range("B27").Value = "dog/dog023.jpg"
$ws.Range("C27").Value = "raten"
$ws.Range("D27").Value = "dog"

$ws.Range("B28").Value = "dog/dog021.jpg"
$ws.Range("C28").Value = "tauschen"
$ws.Range("D28").Value = "dog"

$ws.Range("B29").Value = "dog/dog014.jpg"
$ws.Range("C29").Value = "backen"
$ws.Range("D29").Value = "dog"

$ws.Range("B30").Value = "dog/dog003.jpg"
$ws.Range("C30").Value = "laufen"
$ws.Range("D30").Value = "dog"

$ws.Range("B31").Value = "dog/dog027.jpg"
$ws.Range("C31").Value = "bleiben"
$ws.Range("D31").Value = "dog"

$ws.Range("B32").Value = "house/house000.jpg"
$ws.Range("C32").Value = "fesseln"
$ws.Range("D32").Value = "house"

$ws.Range("B33").Value = "dog/dog001.jpg"
$ws.Range("C33").Value = "schenken"
$ws.Range("D33").Value = "dog"
